$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the last history-table row added for TFS25490 (10/18/2022 entry).
#    The commit rolls this entry back out of the "Document History" table.
# ---------------------------------------------------------------------------
$histTable = $d.Tables.Item(1)
$lastRowIndex = $histTable.Rows.Count
$lastRow = $histTable.Rows.Item($lastRowIndex)
if ($lastRow.Cells.Item(1).Range.Text -like "10/18/2022*") {
    $lastRow.Delete()
}

# ---------------------------------------------------------------------------
# 2. "250 bytes,  Alpha/Numeric" -> "50 bytes,  Alpha/Numeric" (first cell,
#    next to strSubmitterEmail / eMail address of the submitter). Keep the
#    leading "5" as its own run (it was already a separate run before the
#    edit) while the remainder collapses into a single clean run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("250 bytes,  Alpha/Numeric", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $rng.Start
    $d.Range($start, $start + 1).Delete()
    $rest = $d.Range($start + 1, $start + 24)
    $rest.Text = "PLACEHOLDERXXXXXXXXXXXXX"
    $rest2 = $d.Range($start + 1, $start + 1 + 24)
    $rest2.Text = "0 bytes,  Alpha/Numeric"
    $fiveRange = $d.Range($start, $start + 1)
    $fiveRange.Font.Bold = 1
    $fiveRange.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# 3. "250 bytes,  Alpha/Numeric" -> "50 bytes,  Alpha/Numeric" (second cell,
#    near CCO_Report...). This one collapses fully into one run.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("250 bytes,  Alpha/Numeric", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $start2 = $rng2.Start
    $full2 = $d.Range($start2, $start2 + 25)
    $full2.Text = "50 bytes,  Alpha/Numeric"
}

# ---------------------------------------------------------------------------
# 4. Split the single run " or Opportunity" into " or " + "Opportunity"
#    (same visible text/formatting, just two runs instead of one).
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute(" or Opportunity", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $start3 = $rng3.Start
    $opportunityRange = $d.Range($start3 + 4, $start3 + 15)
    $opportunityRange.Font.Bold = 1
    $opportunityRange.Font.Bold = 0
}
